# Auto-generated Excel COM-interop script to apply Spriggan_Profits.xlsx leve-profit updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 876.1667
$ws.Range("I5").Value = 876.1667
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 876.1667
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -761.1667

$ws.Range("H6").Value = 516.375
$ws.Range("I6").Value = 404.53845
$ws.Range("J6").Value = 1001
$ws.Range("K6").Value = 1213.61535
$ws.Range("L6").Value = 3003
$ws.Range("M6").Value = -1101.61535

$ws.Range("H9").Value = 9994.637000000001
$ws.Range("I9").Value = 16907
$ws.Range("J9").Value = 1699.8
$ws.Range("K9").Value = 16907
$ws.Range("L9").Value = 1699.8
$ws.Range("M9").Value = -16738

$ws.Range("H29").Value = 2175.1333
$ws.Range("I29").Value = 247.42857
$ws.Range("J29").Value = 3861.875
$ws.Range("K29").Value = 742.28571
$ws.Range("L29").Value = 11585.625
$ws.Range("M29").Value = -461.28571
$ws.Range("N29").Value = -12147.625

$ws.Range("H51").Value = 3000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 3000
$ws.Range("N51").Value = -3968

$ws.Range("H98").Value = 1370.6111
$ws.Range("I98").Value = 1418.25
$ws.Range("J98").Value = 989.5
$ws.Range("K98").Value = 1418.25
$ws.Range("L98").Value = 989.5
$ws.Range("M98").Value = 79.75
$ws.Range("N98").Value = -3985.5

$ws.Range("H100").Value = 28751.25
$ws.Range("I100").Value = 28751.25
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 28751.25
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -28210.25
$ws.Range("N100").ClearContents()

$ws.Range("H101").Value = 461.16666
$ws.Range("I101").Value = 393.5
$ws.Range("J101").Value = 495
$ws.Range("K101").Value = 1180.5
$ws.Range("L101").Value = 1485
$ws.Range("M101").Value = 441.5
$ws.Range("N101").Value = -4729

$ws.Range("H113").Value = 2511.6667
$ws.Range("I113").Value = 2511.6667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2511.6667
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 742.3332999999998
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 1370.6111
$ws.Range("I122").Value = 1418.25
$ws.Range("J122").Value = 989.5
$ws.Range("K122").Value = 4254.75
$ws.Range("L122").Value = 2968.5
$ws.Range("M122").Value = -1804.75
$ws.Range("N122").Value = -7868.5

$ws.Range("H135").Value = 45455324
$ws.Range("I135").Value = 45455324
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 409097916
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -409095381

$ws.Range("H138").Value = 2453.56
$ws.Range("I138").Value = 2555.875
$ws.Range("J138").Value = 2405.4119
$ws.Range("K138").Value = 7667.625
$ws.Range("L138").Value = 7216.2357
$ws.Range("M138").Value = -2527.625
$ws.Range("N138").Value = -17496.2357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5782.95
$ws.Range("I32").Value = 3907.9312
$ws.Range("J32").Value = 10726.182
$ws.Range("K32").Value = 3907.9312
$ws.Range("L32").Value = 10726.182
$ws.Range("M32").Value = -3620.9312
$ws.Range("N32").Value = -11300.182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 19234170
$ws.Range("I134").Value = 20835352
$ws.Range("J134").Value = 19999.5
$ws.Range("K134").Value = 62506056
$ws.Range("L134").Value = 59998.5
$ws.Range("M134").Value = -62503521
$ws.Range("N134").Value = -65068.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 2024.75
$ws.Range("I33").Value = 2024.75
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 2024.75
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -1645.75
$ws.Range("N33").ClearContents()

$ws.Range("H95").Value = 11916.333
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 11916.333
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 11916.333
$ws.Range("N95").Value = -17408.333

$ws.Range("H111").Value = 38999.5
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 38999.5
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 38999.5
$ws.Range("N111").Value = -47179.5

$ws.Range("H132").Value = 22728814
$ws.Range("I132").Value = 24391762
$ws.Range("J132").Value = 1863
$ws.Range("K132").Value = 73175286
$ws.Range("L132").Value = 5589
$ws.Range("M132").Value = -73172756

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 300
$ws.Range("I36").Value = 300
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 900
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -731

$ws.Range("H46").Value = 100
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 100
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 300
$ws.Range("N46").Value = -482

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 500
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 500
$ws.Range("N19").Value = -1076
$ws.Range("M19").ClearContents()

$ws.Range("H64").Value = 61600
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 61600
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 61600
$ws.Range("N64").Value = -62096

$ws.Range("H67").Value = 61600
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 61600
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 61600
$ws.Range("N67").Value = -63316

$ws.Range("H92").Value = 5856.3
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 5856.3
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 5856.3
$ws.Range("N92").Value = -9600.299999999999

$ws.Range("H123").Value = 49689.6
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 49689.6
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 49689.6
$ws.Range("N123").Value = -54589.6

$ws.Range("H132").Value = 4171438.5
$ws.Range("I132").Value = 4810699.5
$ws.Range("J132").Value = 16240
$ws.Range("K132").Value = 14432098.5
$ws.Range("L132").Value = 48720
$ws.Range("M132").Value = -14429568.5

$ws.Range("H141").Value = 89900
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 89900
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 89900
$ws.Range("N141").Value = -100260

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 999
$ws.Range("I2").Value = 999
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 999
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -887

$ws.Range("H7").Value = 3899
$ws.Range("I7").Value = 3899
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3899
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -3787

$ws.Range("H11").Value = 7000
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 7000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 7000
$ws.Range("N11").Value = -7280
$ws.Range("M11").ClearContents()

$ws.Range("H16").Value = 1684.8823
$ws.Range("I16").Value = 386.65
$ws.Range("J16").Value = 3539.5
$ws.Range("K16").Value = 386.65
$ws.Range("L16").Value = 3539.5
$ws.Range("M16").Value = -216.65

$ws.Range("H22").Value = 4357
$ws.Range("I22").Value = 4125
$ws.Range("J22").Value = 4666.3335
$ws.Range("K22").Value = 4125
$ws.Range("L22").Value = 4666.3335
$ws.Range("M22").Value = -3830
$ws.Range("N22").Value = -5256.3335

$ws.Range("H27").Value = 4357
$ws.Range("I27").Value = 4125
$ws.Range("J27").Value = 4666.3335
$ws.Range("K27").Value = 4125
$ws.Range("L27").Value = 4666.3335
$ws.Range("M27").Value = -4018
$ws.Range("N27").Value = -4880.3335

$ws.Range("H46").Value = 998.6667
$ws.Range("I46").Value = 998.6667
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 998.6667
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -810.6667

$ws.Range("H126").Value = 3899
$ws.Range("I126").Value = 3899
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 11697
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9227

$ws.Range("H132").Value = 68617290
$ws.Range("I132").Value = 68617290
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 205851870
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -205849340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 10000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 10000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 10000
$ws.Range("N19").Value = -10348

$ws.Range("H62").Value = 6575.5
$ws.Range("I62").Value = 6575.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 6575.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -5951.5

$ws.Range("H65").Value = 6575.5
$ws.Range("I65").Value = 6575.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 32877.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -29757.5

$ws.Range("H96").Value = 1303.8125
$ws.Range("I96").Value = 939.2
$ws.Range("J96").Value = 1469.5454
$ws.Range("K96").Value = 939.2
$ws.Range("L96").Value = 1469.5454
$ws.Range("M96").Value = 433.8
$ws.Range("N96").Value = -4215.5454
